$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing
$ws.Range("H43").Value = 6833.3335
$ws.Range("I43").Value = 5250
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 5250
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -5181
$ws.Range("N43").Value = -10138

# Row 51: A Bile Business
$ws.Range("H51").Value = 14778.909
$ws.Range("J51").Value = 14757.7
$ws.Range("L51").Value = 14757.7
$ws.Range("N51").Value = -15725.7

# Row 80: Cleansing the Wicked Humours
$ws.Range("H80").Value = 492.22223
$ws.Range("I80").Value = 576
$ws.Range("J80").Value = 199
$ws.Range("K80").Value = 1728
$ws.Range("L80").Value = 597
$ws.Range("M80").Value = -730
$ws.Range("N80").Value = -2593

# Row 83: Washing Away the Sins (L)
$ws.Range("H83").Value = 492.22223
$ws.Range("I83").Value = 576
$ws.Range("J83").Value = 199
$ws.Range("K83").Value = 5184
$ws.Range("L83").Value = 1791
$ws.Range("M83").Value = -192
$ws.Range("N83").Value = -11775

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 19429596
$ws.Range("J88").Value = 2755747.5
$ws.Range("L88").Value = 2755747.5
$ws.Range("N88").Value = -2756559.5

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 19429596
$ws.Range("J91").Value = 2755747.5
$ws.Range("L91").Value = 2755747.5
$ws.Range("N91").Value = -2758555.5

# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 574.6070999999999
$ws.Range("I92").Value = 384.5238
$ws.Range("K92").Value = 384.5238
$ws.Range("M92").Value = 863.4762000000001

# Row 112: Making Ends Meet
$ws.Range("H112").Value = 41454.71
$ws.Range("I112").Value = 78585.84
$ws.Range("J112").Value = 29077.666
$ws.Range("K112").Value = 235757.52
$ws.Range("L112").Value = 87232.99800000001
$ws.Range("M112").Value = -234649.52
$ws.Range("N112").Value = -89448.99800000001

# Row 113: Amaro Kart
$ws.Range("H113").Value = 3141.1
$ws.Range("I113").Value = 2950.75
$ws.Range("K113").Value = 2950.75
$ws.Range("M113").Value = 303.25

# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2190.4211
$ws.Range("J132").Value = 2166.5
$ws.Range("L132").Value = 6499.5
$ws.Range("N132").Value = -11559.5

# Row 138: All-night Crafting
$ws.Range("H138").Value = 5210.353
$ws.Range("I138").Value = 2425.2222
$ws.Range("J138").Value = 5807.1665
$ws.Range("K138").Value = 7275.6666
$ws.Range("L138").Value = 17421.4995
$ws.Range("M138").Value = -2135.6666
$ws.Range("N138").Value = -27701.4995

$ws = $wb.Worksheets.Item("ARM")
# Row 30: Not Enough Headroom
$ws.Range("H30").Value = 9999
$ws.Range("I30").Value = 9999
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 9999
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -9849
$ws.Range("N30").ClearContents()

# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2002358
$ws.Range("I45").Value = 2002358
$ws.Range("K45").Value = 2002358
$ws.Range("M45").Value = -2001981

# Row 88: The Mast Chance
$ws.Range("H88").Value = 1816.4286
$ws.Range("J88").Value = 2003
$ws.Range("L88").Value = 2003
$ws.Range("N88").Value = -2815

# Row 91: The Rose and the Riveter (L)
$ws.Range("H91").Value = 1816.4286
$ws.Range("J91").Value = 2003
$ws.Range("L91").Value = 2003
$ws.Range("N91").Value = -4811

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4765920.5
$ws.Range("I132").Value = 5003719
$ws.Range("K132").Value = 15011157
$ws.Range("M132").Value = -15008627

$ws = $wb.Worksheets.Item("BSM")
# Row 10: Bring Me the Head Knife of Al'bedo Derssia
$ws.Range("H10").Value = 1000
$ws.Range("I10").Value = 1000
$ws.Range("K10").Value = 1000
$ws.Range("M10").Value = -860

# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 4483.3335
$ws.Range("I20").Value = 2627.3333
$ws.Range("J20").Value = 7267.3335
$ws.Range("K20").Value = 2627.3333
$ws.Range("L20").Value = 7267.3335
$ws.Range("M20").Value = -2380.3333
$ws.Range("N20").Value = -7761.3335

# Row 88: Swords for Plowshares
$ws.Range("H88").Value = 26360.75
$ws.Range("J88").Value = 26360.75
$ws.Range("L88").Value = 26360.75
$ws.Range("N88").Value = -27172.75

# Row 91: Negative, They Are Meat Popsicles (L)
$ws.Range("H91").Value = 26360.75
$ws.Range("J91").Value = 26360.75
$ws.Range("L91").Value = 26360.75
$ws.Range("N91").Value = -29168.75

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2399.6
$ws.Range("I105").Value = 2399.6
$ws.Range("K105").Value = 2399.6
$ws.Range("M105").Value = -652.5999999999999

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 43752200
$ws.Range("I134").Value = 43752200
$ws.Range("K134").Value = 131256600
$ws.Range("M134").Value = -131254065

$ws = $wb.Worksheets.Item("CRP")
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 125003360
$ws.Range("I132").Value = 125003360
$ws.Range("K132").Value = 375010080
$ws.Range("M132").Value = -375007550

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food
$ws.Range("H2").Value = 127.89474
$ws.Range("I2").Value = 122.84615
$ws.Range("J2").Value = 138.83333
$ws.Range("K2").Value = 737.0769
$ws.Range("L2").Value = 832.9999799999999
$ws.Range("M2").Value = -624.0769
$ws.Range("N2").Value = -1058.99998

# Row 5: What a Sap
$ws.Range("H5").Value = 55865.39
$ws.Range("I5").Value = 55865.39
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 167596.17
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -167484.17
$ws.Range("N5").ClearContents()

# Row 34: Fever Pitch
$ws.Range("H34").Value = 1443.2858
$ws.Range("I34").Value = 102
$ws.Range("J34").Value = 1666.8334
$ws.Range("K34").Value = 306
$ws.Range("L34").Value = 5000.5002
$ws.Range("M34").Value = -222
$ws.Range("N34").Value = -5168.5002

# Row 70: Persona non Gratin
$ws.Range("H70").Value = 8074.294
$ws.Range("I70").Value = 4944.077
$ws.Range("K70").Value = 14832.231
$ws.Range("M70").Value = -14517.231

# Row 73: Recipe for Disaster (L)
$ws.Range("H73").Value = 8074.294
$ws.Range("I73").Value = 4944.077
$ws.Range("K73").Value = 14832.231
$ws.Range("M73").Value = -13740.231

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 1818.4667
$ws.Range("I131").Value = 1352.4546
$ws.Range("J131").Value = 3100
$ws.Range("K131").Value = 4057.3638
$ws.Range("L131").Value = 9300
$ws.Range("M131").Value = 982.6361999999999
$ws.Range("N131").Value = -19380

# Row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 55865.39
$ws.Range("I135").Value = 55865.39
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 502788.51
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -500253.51
$ws.Range("N135").ClearContents()

# Row 137: Creative Chocolate
$ws.Range("H137").Value = 11113936
$ws.Range("J137").Value = 3404.6667
$ws.Range("L137").Value = 10214.0001
$ws.Range("N137").Value = -20414.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4277.154
$ws.Range("I80").Value = 4140.75
$ws.Range("K80").Value = 4140.75
$ws.Range("M80").Value = -3142.75

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4277.154
$ws.Range("I83").Value = 4140.75
$ws.Range("K83").Value = 20703.75
$ws.Range("M83").Value = -15711.75

# Row 92: Play It by Ear
$ws.Range("H92").Value = 10000
$ws.Range("J92").Value = 10000
$ws.Range("L92").Value = 10000
$ws.Range("N92").Value = -13744

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value = 663.4167
$ws.Range("I97").Value = 516.1
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 516.1
$ws.Range("L97").Value = 1400
$ws.Range("M97").Value = -20.10000000000002
$ws.Range("N97").Value = -2392

# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 700
$ws.Range("I102").Value = 700
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 700
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 922
$ws.Range("N102").ClearContents()

# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 504499.5
$ws.Range("I113").Value = 999999
$ws.Range("J113").Value = 9000
$ws.Range("K113").Value = 999999
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -997829
$ws.Range("N113").Value = -13340

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 58775.41
$ws.Range("I122").Value = 77566.25
$ws.Range("J122").Value = 8666.5
$ws.Range("K122").Value = 232698.75
$ws.Range("L122").Value = 25999.5
$ws.Range("M122").Value = -230248.75
$ws.Range("N122").Value = -30899.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 20836858
$ws.Range("I132").Value = 20836858
$ws.Range("K132").Value = 62510574
$ws.Range("M132").Value = -62508044

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 2224.5
$ws.Range("I22").Value = 2543.2727
$ws.Range("J22").Value = 1523.2
$ws.Range("K22").Value = 2543.2727
$ws.Range("L22").Value = 1523.2
$ws.Range("M22").Value = -2248.2727
$ws.Range("N22").Value = -2113.2

# Row 27: Fire and Hide
$ws.Range("H27").Value = 2224.5
$ws.Range("I27").Value = 2543.2727
$ws.Range("J27").Value = 1523.2
$ws.Range("K27").Value = 2543.2727
$ws.Range("L27").Value = 1523.2
$ws.Range("M27").Value = -2436.2727
$ws.Range("N27").Value = -1737.2

# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1031.3334
$ws.Range("I46").Value = 1057.6
$ws.Range("J46").Value = 900
$ws.Range("K46").Value = 1057.6
$ws.Range("L46").Value = 900
$ws.Range("M46").Value = -869.5999999999999
$ws.Range("N46").Value = -1276

# Row 122: Hell on Leather
$ws.Range("H122").Value = 837462.8
$ws.Range("I122").Value = 1115285
$ws.Range("K122").Value = 3345855
$ws.Range("M122").Value = -3343405

$ws = $wb.Worksheets.Item("WVR")
# Row 94: Proper Props
$ws.Range("H94").Value = 44444
$ws.Range("J94").Value = 44444
$ws.Range("L94").Value = 44444
$ws.Range("N94").Value = -46246

# Row 126: A Polished Purchase
$ws.Range("H126").Value = 1163.2069
$ws.Range("I126").Value = 1204.1923
$ws.Range("J126").Value = 808
$ws.Range("K126").Value = 3612.5769
$ws.Range("L126").Value = 2424
$ws.Range("M126").Value = -1142.5769
$ws.Range("N126").Value = -7364

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 14711018
$ws.Range("I132").Value = 19232898
$ws.Range("K132").Value = 57698694
$ws.Range("M132").Value = -57696164

Write-Output "Done applying edits"